{"js": "// Insert five new pinyin-annotated sentences, followed by one trailing\n// blank paragraph, right after the existing \"\u53bb\u94f6\u884c(h\u00e1ng)...\" paragraph\n// and before the document's pre-existing trailing empty paragraph.\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that ends the \"\u53bb\u94f6\u884c...\u865a\u6b64\u884c(x\u00edng)\u3002\" sentence \u2014\n// the anchor after which the new content is inserted.\nconst items = body.paragraphs.items;\nlet anchor = null;\nfor (const p of items) {\n  if (p.text.indexOf(\"\u53bb\u94f6\u884c(h\u00e1ng)\u53d6\u94b1\u771f\u662f\u4e0d(b\u00f9)\u865a\u6b64\u884c(x\u00edng)\u3002\") !== -1) {\n    anchor = p;\n  }\n}\nif (!anchor) {\n  // Fallback: insert right before the final (blank) paragraph.\n  anchor = items.length >= 2 ? items[items.length - 2] : items[items.length - 1];\n}\n\nconst newSentences = [\n  \"\u8fd8(h\u00e1i)\u6709(y\u01d2u)\uff0c\u501f\u4f60\u7684(de)\u4e66\u8be5\u8fd8(h\u00e1i)\u6211\u4e86(le)\u3002\",\n  \"\u4f60\u8bf4(shu\u014d)\u8fc7(gu\u00f2)\u7684(de)\u8bdd\u7684(d\u00ed)\u786e\u6709(y\u01d2u)\u9053\u7406\u3002\",\n  \"\u8fd9(zh\u00e8)\u662f\u4e00\u79cd(zh\u01d2ng)\u79d1\u5b66\u79cd(zh\u00f2ng)\u7530\u7684(de)\u65b9\u6cd5\u3002\",\n  \"\u6211\u4eec\u62bd\u7a7a(k\u014dng)\u53bb\u5317(b\u011bi)\u4eac\u822a\u7a7a(k\u014dng)\u822a\u5929\u5b66\u9662\u73a9\u4e00\u4e0b\u5427(ba)\u3002\",\n  \"\u5e03\u5c14\u4ec0(sh\u00ed)\u7ef4\u514b\u662f\u4ec0(sh\u00e9n)\u4e48(me)\u7ec4\u7ec7(zh\u012b)\",\n];\n\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\n// Build the raw paragraph OOXML: one <w:p> per sentence (same run\n// formatting as the surrounding paragraphs \u2014 \u5b8b\u4f53, 12pt / sz 24), plus a\n// final empty paragraph (no run), matching the diff exactly.\nlet paragraphsXml = \"\";\nfor (const sentence of newSentences) {\n  paragraphsXml +=\n    '<w:p><w:pPr><w:ind w:firstLine=\"420\"/></w:pPr>' +\n    '<w:r><w:rPr><w:rFonts w:ascii=\"\u5b8b\u4f53\" w:hAnsi=\"\u5b8b\u4f53\"/><w:sz w:val=\"24\"/></w:rPr>' +\n    '<w:t xml:space=\"preserve\">' + escapeXml(sentence) + '</w:t></w:r></w:p>';\n}\nparagraphsXml += '<w:p><w:pPr><w:ind w:firstLine=\"420\"/></w:pPr></w:p>';\n\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + paragraphsXml + '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nconst insertionPoint = anchor.getRange(Word.RangeLocation.after);\ninsertionPoint.insertOoxml(ooxml, Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Insert five new pinyin-annotated sentences, followed by one trailing\n# blank paragraph, right after the existing \"\u53bb\u94f6\u884c(h\u00e1ng)...\" paragraph\n# and before the document's pre-existing trailing empty paragraph.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"\u53bb\u94f6\u884c(h\u00e1ng)\u53d6\u94b1\u771f\u662f\u4e0d(b\u00f9)\u865a\u6b64\u884c(x\u00edng)\u3002\"\n\n$sentences = @(\n    \"\u8fd8(h\u00e1i)\u6709(y\u01d2u)\uff0c\u501f\u4f60\u7684(de)\u4e66\u8be5\u8fd8(h\u00e1i)\u6211\u4e86(le)\u3002\",\n    \"\u4f60\u8bf4(shu\u014d)\u8fc7(gu\u00f2)\u7684(de)\u8bdd\u7684(d\u00ed)\u786e\u6709(y\u01d2u)\u9053\u7406\u3002\",\n    \"\u8fd9(zh\u00e8)\u662f\u4e00\u79cd(zh\u01d2ng)\u79d1\u5b66\u79cd(zh\u00f2ng)\u7530\u7684(de)\u65b9\u6cd5\u3002\",\n    \"\u6211\u4eec\u62bd\u7a7a(k\u014dng)\u53bb\u5317(b\u011bi)\u4eac\u822a\u7a7a(k\u014dng)\u822a\u5929\u5b66\u9662\u73a9\u4e00\u4e0b\u5427(ba)\u3002\",\n    \"\u5e03\u5c14\u4ec0(sh\u00ed)\u7ef4\u514b\u662f\u4ec0(sh\u00e9n)\u4e48(me)\u7ec4\u7ec7(zh\u012b)\"\n)\n\n# Build a single replacement that keeps the anchor sentence intact, then\n# appends a paragraph break (\"^p\") before each new sentence, plus one more\n# trailing paragraph break for the new blank paragraph. Using Find/Replace\n# (rather than Range.InsertParagraphAfter) keeps the newly split paragraph\n# marks clean (no stray empty run), matching how the original document's\n# own blank paragraph is structured.\n$replacement = $anchorText + \"^p\" + ($sentences -join \"^p\") + \"^p\"\n\n$rng = $d.Range()\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$found = $find.Execute(\n    $anchorText,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    [Microsoft.Office.Interop.Word.WdFindWrap]::wdFindContinue,\n    $false,\n    $replacement,\n    [Microsoft.Office.Interop.Word.WdReplace]::wdReplaceAll\n)\n\nif ($found) {\n    # Locate the paragraph index of the anchor sentence, then apply the\n    # surrounding run formatting (\u5b8b\u4f53, 12pt / sz 24) to just the new\n    # sentence text -- excluding the trailing blank paragraph and\n    # excluding each paragraph's final mark character -- so the new runs\n    # match the rest of the document exactly.\n    $anchorIndex = 0\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq $anchorText) {\n            $anchorIndex = $i\n        }\n    }\n    if ($anchorIndex -eq 0) {\n        $anchorIndex = 3\n    }\n\n    $firstNewIndex = $anchorIndex + 1\n    $lastNewTextIndex = $anchorIndex + $sentences.Count\n\n    $startPos = $d.Paragraphs.Item($firstNewIndex).Range.Start\n    $endPos = $d.Paragraphs.Item($lastNewTextIndex).Range.End - 1  # exclude paragraph mark\n\n    $textRange = $d.Range($startPos, $endPos)\n    $textRange.Font.Name = \"\u5b8b\u4f53\"\n    $textRange.Font.Size = 12\n}\n\n\"paragraphs=\" + $d.Paragraphs.Count\n"}
